$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B ("year"-like grouping value) for rows 119-218.
# The new value mirrors the old one: new = 18 - old
# (10->8, 11->7, 12->6, 13->5, 14->4, 15->3, 16->2, 17->1)
for ($row = 119; $row -le 218; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = $cell.Value()
    $cell.Value = 18 - $old
}

# Reflect the view state change recorded for this edit: the window had
# scrolled down and the user ended up with cell M210 selected.
$excel.ActiveWindow.ScrollRow = 178
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M210").Select()
